$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1: copy format from E1 (bold header style), then set text
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Data cells F2:F52 - time_taken timestamps (no special style, matches column E data cells)
$timeTaken = @(
    "2021-10-05 10:52:12.050061",
    "2021-10-05 10:52:12.050074",
    "2021-10-05 10:52:12.050078",
    "2021-10-05 10:52:12.050082",
    "2021-10-05 10:52:12.050085",
    "2021-10-05 10:52:12.050088",
    "2021-10-05 10:52:12.050091",
    "2021-10-05 10:52:12.050094",
    "2021-10-05 10:52:12.050098",
    "2021-10-05 10:52:12.050101",
    "2021-10-05 10:52:12.050104",
    "2021-10-05 10:52:12.050107",
    "2021-10-05 10:52:12.050110",
    "2021-10-05 10:52:12.050113",
    "2021-10-05 10:52:12.050116",
    "2021-10-05 10:52:12.050120",
    "2021-10-05 10:52:12.050123",
    "2021-10-05 10:52:12.050126",
    "2021-10-05 10:52:12.050129",
    "2021-10-05 10:52:12.050132",
    "2021-10-05 10:52:12.050135",
    "2021-10-05 10:52:12.050138",
    "2021-10-05 10:52:12.050141",
    "2021-10-05 10:52:12.050144",
    "2021-10-05 10:52:12.050148",
    "2021-10-05 10:52:12.050151",
    "2021-10-05 10:52:12.050154",
    "2021-10-05 10:52:12.050157",
    "2021-10-05 10:52:12.050160",
    "2021-10-05 10:52:12.050163",
    "2021-10-05 10:52:12.050166",
    "2021-10-05 10:52:12.050169",
    "2021-10-05 10:52:12.050173",
    "2021-10-05 10:52:12.050176",
    "2021-10-05 10:52:12.050179",
    "2021-10-05 10:52:12.050182",
    "2021-10-05 10:52:12.050185",
    "2021-10-05 10:52:12.050188",
    "2021-10-05 10:52:12.050192",
    "2021-10-05 10:52:12.050195",
    "2021-10-05 10:52:12.050199",
    "2021-10-05 10:52:12.050225",
    "2021-10-05 10:52:12.050228",
    "2021-10-05 10:52:12.050248",
    "2021-10-05 10:52:12.050251",
    "2021-10-05 10:52:12.050254",
    "2021-10-05 10:52:12.050257",
    "2021-10-05 10:52:12.050290",
    "2021-10-05 10:52:12.050316",
    "2021-10-05 10:52:12.050324",
    "2021-10-05 10:52:12.050327"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}

Write-Host "Added time_taken column F with header and $($timeTaken.Length) data rows"
